$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Update "Status" text (used by Overview!E2:F3 and the Status column on the
#    per-locale sheets) from "Ready for handoff" to the handed-back message.
# ---------------------------------------------------------------------------
$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$ws1.Cells.Replace($oldStatus, $newStatus, 1) | Out-Null
$ws2.Cells.Replace($oldStatus, $newStatus, 1) | Out-Null
$ws3.Cells.Replace($oldStatus, $newStatus, 1) | Out-Null

# ---------------------------------------------------------------------------
# 2. Populate the "Latest Target File" / "Latest Handback File" (and, for
#    de-de, "Latest Handback DateTime") columns for the two locale sheets,
#    which is what "generating the handback report" actually produces.
# ---------------------------------------------------------------------------

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3fea4bd994ffd600339ea067d7e22e913f85a31/e2e/594b5d25-a9d7-4770-b5d6-a16a85dbd574.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f3fea4bd994ffd600339ea067d7e22e913f85a31/e2e/c87c3933-ca4b-430a-b874-d5c96a0cb62e.md"

$mdName1 = "594b5d25-a9d7-4770-b5d6-a16a85dbd574.md"
$mdName2 = "c87c3933-ca4b-430a-b874-d5c96a0cb62e.md"

# -- zh-cn sheet --------------------------------------------------------
$ws2.Range("J2").Value = "594b5d25-a9d7-4770-b5d6-a16a85dbd574.d19b2f296142dda1fbfdf418416143531bbd21e5.zh-cn.xlf"
$ws2.Range("J3").Value = "c87c3933-ca4b-430a-b874-d5c96a0cb62e.22ff6ed4ab21fdb63eb6b2d237911a0ef19d6585.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-07 09:41:11"
$ws2.Range("K3").Value = "2016-09-07 09:41:11"

$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null

# -- de-de sheet ----------------------------------------------------------
$ws3.Range("J2").Value = "594b5d25-a9d7-4770-b5d6-a16a85dbd574.d19b2f296142dda1fbfdf418416143531bbd21e5.de-de.xlf"
$ws3.Range("J3").Value = "c87c3933-ca4b-430a-b874-d5c96a0cb62e.22ff6ed4ab21fdb63eb6b2d237911a0ef19d6585.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-07 09:41:36"
$ws3.Range("K3").Value = "2016-09-07 09:41:36"

$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Widen the columns that now hold longer text (the Status column as well
#    as the newly filled in "Latest Target File"/"Latest Handback File"
#    columns).
# ---------------------------------------------------------------------------
function Set-ExactColumnWidth($col, [double]$targetOoxmlWidth) {
    $col.ColumnWidth = $targetOoxmlWidth - (5.0/6.0)
}

Set-ExactColumnWidth $ws1.Columns.Item(5) 29.9777050018311   # Overview!E (zh-cn status)
Set-ExactColumnWidth $ws1.Columns.Item(6) 29.9777050018311   # Overview!F (de-de status)

Set-ExactColumnWidth $ws2.Columns.Item(3) 29.9777050018311   # zh-cn!C (Status)
Set-ExactColumnWidth $ws2.Columns.Item(9) 40                 # zh-cn!I (Latest Target File)
Set-ExactColumnWidth $ws2.Columns.Item(10) 40                # zh-cn!J (Latest Handback File)

Set-ExactColumnWidth $ws3.Columns.Item(3) 29.9777050018311   # de-de!C (Status)
Set-ExactColumnWidth $ws3.Columns.Item(9) 40                 # de-de!I (Latest Target File)
Set-ExactColumnWidth $ws3.Columns.Item(10) 40                # de-de!J (Latest Handback File)

Write-Host "Generated handback report"
